$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "58.588.93"
$ws.Cells.Item(2,5).Value = "  -1.76%  "

$ws.Cells.Item(3,4).Value = "2.615.88"
$ws.Cells.Item(3,5).Value = "  -0.02%  "

$ws.Cells.Item(4,5).Value = "  -0.03%  "

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "535.93"
$ws.Cells.Item(5,4).ClearFormats()
$ws.Cells.Item(5,5).Value = "  -0.43%  "

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "143.14"
$ws.Cells.Item(6,4).ClearFormats()
$ws.Cells.Item(6,5).Value = "  +0.54%  "

$ws.Cells.Item(8,5).Value = "  +3.41%  "

$ws.Cells.Item(9,4).Value = "2.620.02"
$ws.Cells.Item(9,5).Value = "  -0.12%  "

$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "6.84"
$ws.Cells.Item(10,4).ClearFormats()
$ws.Cells.Item(10,5).Value = "  +4.11%  "

$ws.Cells.Item(11,5).Value = "  -1.90%  "

$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "0.334"
$ws.Cells.Item(12,4).ClearFormats()
$ws.Cells.Item(12,5).Value = "  -0.43%  "

$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = "0.137"
$ws.Cells.Item(13,4).ClearFormats()
$ws.Cells.Item(13,5).Value = "  +2.22%  "

$ws.Cells.Item(14,4).Value = "3.070.40"

$ws.Cells.Item(15,4).Value = "58.539.02"
$ws.Cells.Item(15,5).Value = "  -1.72%  "

$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "20.83"
$ws.Cells.Item(16,4).ClearFormats()
$ws.Cells.Item(16,5).Value = "  +0.27%  "

$ws.Cells.Item(17,4).Value = "2.604.39"
$ws.Cells.Item(17,5).Value = "  -0.42%  "

$ws.Cells.Item(18,5).Value = "  -1.16%  "

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "4.42"
$ws.Cells.Item(19,4).ClearFormats()
$ws.Cells.Item(19,5).Value = "  +1.20%  "

$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "335.44"
$ws.Cells.Item(20,4).ClearFormats()
$ws.Cells.Item(20,5).Value = "  -1.56%  "

$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "10.14"
$ws.Cells.Item(21,4).ClearFormats()
$ws.Cells.Item(21,5).Value = "  +0.14%  "

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "6.21"
$ws.Cells.Item(22,4).ClearFormats()
$ws.Cells.Item(22,5).Value = "  -2.31%  "

$ws.Cells.Item(23,5).Value = "  -0.10%  "

$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "67.06"
$ws.Cells.Item(24,4).ClearFormats()
$ws.Cells.Item(24,5).Value = "  -0.28%  "

$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "0.422"
$ws.Cells.Item(25,4).ClearFormats()
$ws.Cells.Item(25,5).Value = "  +3.01%  "

$ws.Cells.Item(26,5).Value = "  -0.03%  "

$ws.Cells.Item(27,5).Value = "  -2.95%  "

$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = "7.10"
$ws.Cells.Item(28,4).ClearFormats()
$ws.Cells.Item(28,5).Value = "  -1.93%  "

$ws.Cells.Item(29,4).Value = "0.0₃0736"
$ws.Cells.Item(29,5).Value = "  -1.45%  "

$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = "1.65"
$ws.Cells.Item(31,4).ClearFormats()
$ws.Cells.Item(31,5).Value = "  -1.21%  "

$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = "5.93"
$ws.Cells.Item(32,4).ClearFormats()
$ws.Cells.Item(32,5).Value = "  +1.54%  "

$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "153.14"
$ws.Cells.Item(33,4).ClearFormats()
$ws.Cells.Item(33,5).Value = "  +1.57%  "

$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "18.91"
$ws.Cells.Item(34,4).ClearFormats()
$ws.Cells.Item(34,5).Value = "  +0.42%  "

$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = "3.91"
$ws.Cells.Item(35,4).ClearFormats()
$ws.Cells.Item(35,5).Value = "  -2.05%  "

$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "37.10"
$ws.Cells.Item(36,4).ClearFormats()
$ws.Cells.Item(36,5).Value = "  -1.02%  "

$ws.Cells.Item(37,5).Value = "  -1.41%  "

$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = "0.838"
$ws.Cells.Item(38,4).ClearFormats()
$ws.Cells.Item(38,5).Value = "  +0.38%  "

$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "0.827"
$ws.Cells.Item(39,4).ClearFormats()
$ws.Cells.Item(39,5).Value = "  -0.31%  "

$ws.Cells.Item(40,5).Value = "  -2.67%  "

$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "3.60"
$ws.Cells.Item(41,4).ClearFormats()
$ws.Cells.Item(41,5).Value = "  +1.60%  "

$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "283.92"
$ws.Cells.Item(42,4).ClearFormats()
$ws.Cells.Item(42,5).Value = "  +1.89%  "

$ws.Cells.Item(43,5).Value = "  +0.06%  "

$ws.Cells.Item(44,5).Value = "  -1.15%  "

$ws.Cells.Item(45,5).Value = "  -0.37%  "

$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "0.0951"
$ws.Cells.Item(46,4).ClearFormats()
$ws.Cells.Item(46,5).Value = "  +0.10%  "

$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "19.09"
$ws.Cells.Item(47,4).ClearFormats()
$ws.Cells.Item(47,5).Value = "  +1.67%  "

$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = "0.0531"
$ws.Cells.Item(48,4).ClearFormats()
$ws.Cells.Item(48,5).Value = "  +1.20%  "

$ws.Cells.Item(49,5).Value = "  +1.25%  "

$ws.Cells.Item(50,4).Value = "1.942.03"
$ws.Cells.Item(50,5).Value = "  -0.61%  "

$ws.Cells.Item(51,5).Value = "  -1.08%  "
